$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B/C are plain text (coin name / URL) -- direct assignment is safe.
# Column D holds numeric-looking strings formatted as text in the source data
# (e.g. "1.00", "0.999", thousand-dot-separated values). Excel's COM layer
# auto-coerces plain numeric-looking strings into real numbers, which would
# lose the original text formatting (trailing zeros, grouping dots, etc.), so
# we force a quote-prefixed text entry for those cells.
# Column E values contain spaces/percent signs and are kept as text naturally.

$ws.Range("D2").Value = "'49.150.79"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "'2.630.15"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'111.74"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").Value = "'322.79"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.542"
$ws.Range("E9").Value = "  -3.08%  "
$ws.Range("D10").Value = "'39.71"
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("D11").Value = "'19.73"
$ws.Range("E11").Value = "  -5.17%  "
$ws.Range("D12").Value = "'0.0811"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "'7.25"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "'3.036.52"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "'2.630.29"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "'0.859"
$ws.Range("E17").Value = "  -1.97%  "
$ws.Range("D18").Value = "'49.082.51"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  -2.93%  "
$ws.Range("D20").Value = "'12.94"
$ws.Range("E20").Value = "  -3.47%  "
$ws.Range("D21").Value = "'6.69"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").Value = "'0.0₃0945"
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").Value = "'269.31"
$ws.Range("E23").Value = "  -4.42%  "
$ws.Range("D24").Value = "'68.54"
$ws.Range("E24").Value = "  -5.88%  "
$ws.Range("D25").Value = "'2.54"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("D26").Value = "'26.16"
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "'10.18"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").Value = "'35.11"
$ws.Range("E30").Value = "  -3.37%  "
$ws.Range("E31").Value = "  -3.93%  "
$ws.Range("D32").Value = "'49.43"
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0801"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").Value = "'19.03"
$ws.Range("E36").Value = "  -4.07%  "
$ws.Range("D37").Value = "'4.96"
$ws.Range("E37").Value = "  +4.22%  "
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("D40").Value = "'126.50"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").Value = "'0.111"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").Value = "'22.19"
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("E43").Value = "  -4.21%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "'2.067.33"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  +6.49%  "
$ws.Range("E47").Value = "  -4.93%  "
$ws.Range("E48").Value = "  -4.97%  "
$ws.Range("D49").Value = "'8.90"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").Value = "'5.20"
$ws.Range("E50").Value = "  -3.28%  "
$ws.Range("D51").Value = "'58.72"
$ws.Range("E51").Value = "  +1.39%  "
